# Automatische test-sync: 2025-06-19 21:47:50
# Adds the new "Afmelding nieuwsbrief" mail log entry as row 30 on the
# "Logs" sheet, extends the conditional formatting ranges to include it,
# and refreshes the category counts/order on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new incoming mail as row 30 ---------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A30").Value = "Afmelding nieuwsbrief"
$logs.Range("B30").Value = "mailmind.test@zohomail.eu"
$logs.Range("C30").Value = "Graag afmelden voor de nieuwsbrief. Dank u."
$logs.Range("D30").Value = "Afmelding / Nieuwsbrief"
$logs.Range("F30").Value = "2025-06-19 21:47:11"
$logs.Range("G30").Value = "Nee"

# Extend the conditional formatting that highlights Categorie/Beantwoord
# so the new row is covered too (D2:D29 -> D2:D30, G2:G29 -> G2:G30).
$categorieRules = $logs.Range("D2:D29").FormatConditions
for ($i = 1; $i -le $categorieRules.Count(); $i++) {
    $categorieRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D30"))
}

$beantwoordRules = $logs.Range("G2:G29").FormatConditions
for ($i = 1; $i -le $beantwoordRules.Count(); $i++) {
    $beantwoordRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G30"))
}

# --- Dashboard sheet: "Afmelding / Nieuwsbrief" count goes 3 -> 4, ----
# which moves it above "Factuur / Administratie" in the ranking table.
$dashboard = $wb.Worksheets.Item("Dashboard")

$dashboard.Range("A4").Value = "Afmelding / Nieuwsbrief"
$dashboard.Range("B4").Value = 4
$dashboard.Range("A6").Value = "Factuur / Administratie"
$dashboard.Range("B6").Value = 3
